$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row labels: the "International levy on shipping..." row and the
# "NCQG: Developing countries..." row swap places (A4 <-> A5), and the
# full dataset was refreshed with the final computed means.
$ws.Range("A4").Value = "NCQG: Developing countries providing `$300 bn a`nyear in climate finance for developing countries"
$ws.Range("A5").Value = "International levy on shipping carbon emissions,`nreturned to countries based on population"
# Re-assert the default (non-wrapped) row height: writing a multi-line
# string otherwise leaves the row auto-sized/taller than the original.
$ws.Range("A4").EntireRow.AutoFit()
$ws.Range("A5").EntireRow.AutoFit()

# --- Row 2 (Bridgetown initiative) updated values
$ws.Range("B2").Value = 0.53880897005617
$ws.Range("C2").Value = 0.622717294167845
$ws.Range("D2").Value = 0.500531888330043
$ws.Range("E2").Value = 0.569056364627645
$ws.Range("F2").Value = 0.864562234521199
$ws.Range("G2").Value = 0.306050220846904
$ws.Range("H2").Value = 0.685511127034902
$ws.Range("I2").Value = 0.729708508139871
$ws.Range("J2").Value = 0.565384609904008
$ws.Range("K2").Value = 0.381834133583139
$ws.Range("L2").Value = 0.622904912769984
$ws.Range("M2").Value = 0.834571283823028
$ws.Range("N2").Value = 0.455071120277596

# --- Row 3 (L&D) updated values
$ws.Range("B3").Value = 0.457646793262474
$ws.Range("C3").Value = 0.487373574197681
$ws.Range("D3").Value = 0.36591653698727
$ws.Range("E3").Value = 0.434921823697749
$ws.Range("F3").Value = 0.713663570537277
$ws.Range("G3").Value = 0.344465636120429
$ws.Range("H3").Value = 0.650280596850207
$ws.Range("I3").Value = 0.468581077682354
$ws.Range("J3").Value = 0.354161380059063
$ws.Range("K3").Value = 0.278997853200917
$ws.Range("L3").Value = 0.695125891930721
$ws.Range("M3").Value = 0.958227224438065
$ws.Range("N3").Value = 0.364623737617183

# --- Row 4 (now NCQG) updated values
$ws.Range("B4").Value = 0.320103458109216
$ws.Range("C4").Value = 0.385513480941771
$ws.Range("D4").Value = 0.282328612443547
$ws.Range("E4").Value = 0.386595005935138
$ws.Range("F4").Value = 0.564560966610589
$ws.Range("G4").Value = 0.195137593917434
$ws.Range("H4").Value = 0.547858472580792
$ws.Range("I4").Value = 0.355483061828207
$ws.Range("J4").Value = 0.268682896761476
$ws.Range("K4").Value = 0.0605446304055215
$ws.Range("L4").Value = 0.697267347192898
$ws.Range("M4").Value = 0.776081506237869
$ws.Range("N4").Value = 0.168711740455824

# --- Row 5 (now International levy on shipping) updated values
$ws.Range("B5").Value = 0.31788661467437
$ws.Range("C5").Value = 0.420490345740477
$ws.Range("D5").Value = 0.536683844982711
$ws.Range("E5").Value = 0.307558555086059
$ws.Range("F5").Value = 0.574957116219404
$ws.Range("G5").Value = 0.119901621697207
$ws.Range("H5").Value = 0.465856915676184
$ws.Range("I5").Value = 0.441669722225902
$ws.Range("J5").Value = 0.404635709261012
$ws.Range("K5").Value = 0.0564462476899513
$ws.Range("L5").Value = 0.36590464407975
$ws.Range("M5").Value = 0.660342805051329
$ws.Range("N5").Value = 0.262398605005224

# --- Row 6 (International levy on aviation) updated values
$ws.Range("B6").Value = 0.00973878914166017
$ws.Range("C6").Value = 0.0839217421719147
$ws.Range("D6").Value = 0.187669811959432
$ws.Range("E6").Value = 0.0887612365922009
$ws.Range("F6").Value = 0.0801259627100831
$ws.Range("G6").Value = -0.0183594094253829
$ws.Range("H6").Value = 0.0679321304395113
$ws.Range("I6").Value = 0.0779461502618107
$ws.Range("J6").Value = 0.010117471105896
$ws.Range("K6").Value = -0.0959998598824107
$ws.Range("L6").Value = -0.0180680286245144
$ws.Range("M6").Value = 0.414946008655522
$ws.Range("N6").Value = -0.0465102463465209
